$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) holds values that look numeric (e.g. "214.75",
# "1.129.17", "0.0930") but must stay literal text: some use a second "."
# as a thousands separator, and some rely on exact trailing zeros. If we
# just assign .Value, Excel will happily parse the numeric-looking ones as
# floating point numbers and mangle them (e.g. "0.530" -> 0.53). Force
# every Price cell that is about to change to Text format first so the
# assignment below is stored verbatim.
$priceRefs = @(
    "D2","D3","D5","D8","D10","D12","D13","D14","D15","D16","D17","D19",
    "D21","D24","D27","D36","D38","D40","D41","D42","D44","D48","D50"
)
foreach ($ref in $priceRefs) {
    $ws.Range($ref).NumberFormat = "@"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "26.040.10"
$ws.Range("E2").Value = "  +0.49%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.637.61"
$ws.Range("E3").Value = "  +0.05%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.56%  "

# Row 5 - BNB
$ws.Range("D5").Value = "214.75"
$ws.Range("E5").Value = "  -0.31%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  -0.08%  "

# Row 8 - Cardano
$ws.Range("D8").Value = "0.252"
$ws.Range("E8").Value = "  -1.62%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  -1.40%  "

# Row 10 - Solana
$ws.Range("D10").Value = "18.77"
$ws.Range("E10").Value = "  -4.07%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +0.37%  "

# Row 12 - WrappedEther
$ws.Range("D12").Value = "1.704.36"
$ws.Range("E12").Value = "  +4.09%  "

# Row 13 - Polkadot
$ws.Range("D13").Value = "4.22"
$ws.Range("E13").Value = "  -1.41%  "

# Row 14 - Polygon
$ws.Range("D14").Value = "0.533"
$ws.Range("E14").Value = "  -1.95%  "

# Row 15 / 16 - Litecoin and ShibaInu swap places (ShibaInu now ranks above
# Litecoin), each with refreshed price/volume figures.
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").Value = "0.0₃0749"
$ws.Range("E15").Value = "  -2.04%  "

$ws.Range("B16").Value = "Litecoin"
$ws.Range("C16").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D16").Value = "62.31"
$ws.Range("E16").Value = "  -0.76%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "26.061.88"

# Row 19 - BitcoinCash
$ws.Range("D19").Value = "191.48"
$ws.Range("E19").Value = "  -0.69%  "

# Row 20 - Uniswap
$ws.Range("E20").Value = "  -1.54%  "

# Row 21 - Avalanche
$ws.Range("D21").Value = "9.63"
$ws.Range("E21").Value = "  -2.81%  "

# Row 22 - Chainlink
$ws.Range("E22").Value = "  -1.51%  "

# Row 23 - Stellar
$ws.Range("E23").Value = "  +1.28%  "

# Row 24 - Monero
$ws.Range("D24").Value = "143.87"
$ws.Range("E24").Value = "  +0.40%  "

# Row 25 - BinanceUSD
$ws.Range("E25").Value = "  +0.95%  "

# Row 26 - Toncoin
$ws.Range("E26").Value = "  -0.85%  "

# Row 27 - Cosmos
$ws.Range("D27").Value = "6.79"
$ws.Range("E27").Value = "  -1.37%  "

# Row 28 - EthereumClassic
$ws.Range("E28").Value = "  -1.78%  "

# Row 29 - PancakeSwap
$ws.Range("E29").Value = "  -0.33%  "

# Row 30 - Hedera
$ws.Range("E30").Value = "  -2.90%  "

# Row 31 - Filecoin
$ws.Range("E31").Value = "  -2.26%  "

# Row 32 - InternetComputer(DFINITY)
$ws.Range("E32").Value = "  -2.97%  "

# Row 33 - LidoDAOToken
$ws.Range("E33").Value = "  -1.46%  "

# Row 34 - HuobiToken
$ws.Range("E34").Value = "  -0.77%  "

# Row 35 - ARBITRUM
$ws.Range("E35").Value = "  -2.33%  "

# Row 36 - Maker
$ws.Range("D36").Value = "1.129.17"
$ws.Range("E36").Value = "  -0.40%  "

# Row 37 - MXToken
$ws.Range("E37").Value = "  +0.05%  "

# Row 38 - ImmutableX
$ws.Range("D38").Value = "0.530"
$ws.Range("E38").Value = "  -2.43%  "

# Row 39 - VeChain
$ws.Range("E39").Value = "  -0.81%  "

# Row 40 - Quant
$ws.Range("D40").Value = "99.02"
$ws.Range("E40").Value = "  -0.25%  "

# Row 41 - TrustWalletToken
$ws.Range("D41").Value = "0.789"
$ws.Range("E41").Value = "  -1.07%  "

# Row 42 - FraxShare
$ws.Range("D42").Value = "5.31"
$ws.Range("E42").Value = "  -3.03%  "

# Row 43 - BabyDogeCoin
$ws.Range("E43").Value = "  -0.61%  "

# Row 44 - Aave
$ws.Range("D44").Value = "55.63"
$ws.Range("E44").Value = "  -1.73%  "

# Row 45 - Cronos
$ws.Range("E45").Value = "  -0.64%  "

# Row 46 - RenderToken
$ws.Range("E46").Value = "  +1.44%  "

# Row 47 - Mantle
$ws.Range("E47").Value = "  +0.09%  "

# Row 48 - EnergySwap
$ws.Range("D48").Value = "7.62"
$ws.Range("E48").Value = "  -0.99%  "

# Row 49 - USDD
$ws.Range("E49").Value = "  +0.34%  "

# Row 50 - Algorand
$ws.Range("D50").Value = "0.0930"
$ws.Range("E50").Value = "  -2.89%  "

# Row 51 - NEARProtocol
$ws.Range("E51").Value = "  -0.26%  "
